$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save off the original (pre-edit) values for rows 2, 3 and 4 before overwriting
# any of them, since the new row values are derived from a cyclic shift of the
# existing rows' data (row2 <- old row3, row3 <- old row4, row4 <- old row2).
$origD2 = $ws.Range("D2").Value2
$origJ2 = $ws.Range("J2").Value2

$origD3 = $ws.Range("D3").Value2
$origJ3 = $ws.Range("J3").Value2
$origK3 = $ws.Range("K3").Value2
$origL3 = $ws.Range("L3").Value2
$origM3 = $ws.Range("M3").Value2
$origN3 = $ws.Range("N3").Value2
$origO3 = $ws.Range("O3").Value2
$origP3 = $ws.Range("P3").Value2
$origQ3 = $ws.Range("Q3").Value2

$origD4 = $ws.Range("D4").Value2
$origJ4 = $ws.Range("J4").Value2
$origK4 = $ws.Range("K4").Value2
$origL4 = $ws.Range("L4").Value2
$origM4 = $ws.Range("M4").Value2
$origN4 = $ws.Range("N4").Value2
$origO4 = $ws.Range("O4").Value2
$origP4 = $ws.Range("P4").Value2
$origQ4 = $ws.Range("Q4").Value2

# Row 2 gets the date/volume coming from old row 3
$ws.Range("D2").Value = $origD3
$ws.Range("J2").Value = $origJ3

# Row 3 gets old row 4's data
$ws.Range("D3").Value = $origD4
$ws.Range("J3").Value = $origJ4
$ws.Range("K3").Value = $origK4
$ws.Range("L3").Value = $origL4
$ws.Range("M3").Value = $origM4
$ws.Range("N3").Value = $origN4
$ws.Range("O3").Value = $origO4
$ws.Range("P3").Value = $origP4
$ws.Range("Q3").Value = $origQ4

# Row 4 gets old row 2's data
$ws.Range("D4").Value = $origD2
$ws.Range("J4").Value = $origJ2
$ws.Range("K4").Value = $origK3
$ws.Range("L4").Value = $origL3
$ws.Range("M4").Value = $origM3
$ws.Range("N4").Value = $origN3
$ws.Range("O4").Value = $origO3
$ws.Range("P4").Value = $origP3
$ws.Range("Q4").Value = $origQ3
